# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps recorded for the most recent handback run.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet, row 2 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-17 14:53:10"
$wsZh.Range("H2").Value = "2016-03-17 14:53:33"

# --- de-de sheet, row 2 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-17 14:53:13"
$wsDe.Range("H2").Value = "2016-03-17 14:53:39"
